$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 516
$ws1.Range("F3").Value = 6162
$ws1.Range("F4").Value = 391
$ws1.Range("F5").Value = 87
$ws1.Range("F6").Value = 120
$ws1.Range("F7").Value = 5
$ws1.Range("F8").Value = 70
$ws1.Range("F9").Value = 561
$ws1.Range("F10").Value = 42

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 6
$ws2.Range("F3").Value = 8

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 516
$ws4.Range("F3").Value = 6162
$ws4.Range("F4").Value = 391
$ws4.Range("F5").Value = 6
$ws4.Range("F6").Value = 87
$ws4.Range("F7").Value = 120
$ws4.Range("F8").Value = 5
$ws4.Range("F9").Value = 8
$ws4.Range("F10").Value = 70
$ws4.Range("F11").Value = 561
$ws4.Range("F12").Value = 42
